{"js": "// Update the numeric inputs/results and one label in the shear-check\n// tables to reflect a new set of design forces / parameters, per the\n// commit \"Detailed shear para list of forces\".\n//\n// Each table on the page is addressed by its 0-based index in\n// context.document.body.tables (Materials, Geometry, Design forces,\n// Limit checks, Shear reinforcement strength, Concrete strength), and\n// each cell inside it by its 0-based (row, col) position. Using\n// position-based addressing (rather than text search-and-replace)\n// avoids ambiguity from values that repeat verbatim elsewhere in the\n// document (e.g. \"15.0\", \"2.12\", \"10.47\", \"\u2714\ufe0f\").\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst edits = [\n  // Materials\n  { t: 0, r: 1, c: 2, newVal: \"101\" },        // Section Label\n  { t: 0, r: 2, c: 2, newVal: \"25\" },          // Concrete strength (fc)\n  { t: 0, r: 3, c: 2, newVal: \"420\" },         // Steel yield strength (fy)\n  { t: 0, r: 4, c: 2, newVal: \"2500.0\" },      // Concrete density\n  // Geometry\n  { t: 1, r: 1, c: 2, newVal: \"50.0\" },        // Section height (h)\n  { t: 1, r: 2, c: 2, newVal: \"20.0\" },        // Section width (b)\n  { t: 1, r: 3, c: 2, newVal: \"2.5\" },         // Clear cover (cc)\n  { t: 1, r: 4, c: 2, newVal: \"0.0\" },         // Longitudinal tension rebar (As)\n  // Design forces\n  { t: 2, r: 1, c: 2, newVal: \"0\" },           // Axial, positive compression (Nu)\n  { t: 2, r: 2, c: 2, newVal: \"100\" },         // Shear (Vu)\n  // Limit checks\n  { t: 3, r: 1, c: 2, newVal: \"20.0\" },        // Stirrup spacing along length - Value\n  { t: 3, r: 1, c: 4, newVal: \"23.05\" },       // Stirrup spacing along length - Max.\n  { t: 3, r: 2, c: 2, newVal: \"14.4\" },        // Stirrup spacing along width - Value\n  { t: 3, r: 2, c: 4, newVal: \"46.1\" },        // Stirrup spacing along width - Max.\n  { t: 3, r: 3, c: 2, newVal: \"2.83\" },        // Minimum shear reinforcement - Value\n  { t: 3, r: 3, c: 3, newVal: \"1.75\" },        // Minimum shear reinforcement - Min.\n  // Shear reinforcement strength\n  { t: 4, r: 2, c: 2, newVal: \"6.0\" },         // Stirrup diameter (db)\n  { t: 4, r: 3, c: 2, newVal: \"20.0\" },        // Stirrup spacing (s)\n  { t: 4, r: 4, c: 2, newVal: \"46.1\" },        // Effective height (d)\n  { t: 4, r: 5, c: 2, newVal: \"1.75\" },        // Minimum shear reinforcing (Av,min)\n  { t: 4, r: 6, c: 2, newVal: \"2.98\" },        // Required shear reinforcing (Av,req)\n  { t: 4, r: 7, c: 0, newVal: \"Defined shear reinforcing\" }, // label\n  { t: 4, r: 7, c: 2, newVal: \"2.83\" },        // Shear reinforcing (Av)\n  { t: 4, r: 8, c: 2, newVal: \"39.1\" },        // Shear steel strength (\u00d8Vs)\n  // Concrete strength\n  { t: 5, r: 1, c: 2, newVal: \"922.0\" },       // Effective shear area (Acv)\n  { t: 5, r: 2, c: 2, newVal: \"0.0\" },         // Longitudinal reinforcement ratio (\u03c1w)\n  { t: 5, r: 3, c: 2, newVal: \"0.839\" },       // Size modification factor (\u03bbs)\n  { t: 5, r: 4, c: 2, newVal: \"0.0\" },         // Axial stress (\u03c3Nu)\n  { t: 5, r: 5, c: 2, newVal: \"0.85\" },        // Concrete effective shear stress (kc)\n  { t: 5, r: 6, c: 2, newVal: \"58.78\" },       // Concrete strength (\u00d8Vc)\n  { t: 5, r: 7, c: 2, newVal: \"286.97\" },      // Maximum shear strength (\u00d8Vmax)\n  { t: 5, r: 8, c: 2, newVal: \"97.88\" },       // Total shear strength (\u00d8Vn)\n  { t: 5, r: 10, c: 2, newVal: \"1.02\" },       // Demand Capacity Ratio (DCR)\n  { t: 5, r: 10, c: 3, newVal: \"\u274c\" },         // Demand Capacity Ratio - Ok?\n];\n\nfor (const e of edits) {\n  const cell = tables.items[e.t].getCell(e.r, e.c);\n  cell.value = e.newVal;\n}\n\nawait context.sync();\n", "ps1": "# Update the numeric inputs/results and one label in the shear-check\n# tables to reflect a new set of design forces / parameters, per the\n# commit \"Detailed shear para list of forces\".\n#\n# Each table on the page is addressed by its 1-based index in\n# $d.Tables (Materials, Geometry, Design forces, Limit checks, Shear\n# reinforcement strength, Concrete strength), and each cell inside it\n# by its 1-based (row, col) position via Table.Cell(row, col). Using\n# position-based addressing (rather than text search-and-replace)\n# avoids ambiguity from values that repeat verbatim elsewhere in the\n# document (e.g. \"15.0\", \"2.12\", \"10.47\", \"\u2714\ufe0f\").\n\n$d = $word.ActiveDocument\n\nfunction Set-CellText($table, $row, $col, $text) {\n    $table.Cell($row, $col).Range.Text = $text\n}\n\n$tbl1 = $d.Tables.Item(1)   # Materials\nSet-CellText $tbl1 2 3 \"101\"        # Section Label\nSet-CellText $tbl1 3 3 \"25\"         # Concrete strength (fc)\nSet-CellText $tbl1 4 3 \"420\"        # Steel yield strength (fy)\nSet-CellText $tbl1 5 3 \"2500.0\"     # Concrete density\n\n$tbl2 = $d.Tables.Item(2)   # Geometry\nSet-CellText $tbl2 2 3 \"50.0\"       # Section height (h)\nSet-CellText $tbl2 3 3 \"20.0\"       # Section width (b)\nSet-CellText $tbl2 4 3 \"2.5\"        # Clear cover (cc)\nSet-CellText $tbl2 5 3 \"0.0\"        # Longitudinal tension rebar (As)\n\n$tbl3 = $d.Tables.Item(3)   # Design forces\nSet-CellText $tbl3 2 3 \"0\"          # Axial, positive compression (Nu)\nSet-CellText $tbl3 3 3 \"100\"        # Shear (Vu)\n\n$tbl4 = $d.Tables.Item(4)   # Limit checks\nSet-CellText $tbl4 2 3 \"20.0\"       # Stirrup spacing along length - Value\nSet-CellText $tbl4 2 5 \"23.05\"      # Stirrup spacing along length - Max.\nSet-CellText $tbl4 3 3 \"14.4\"       # Stirrup spacing along width - Value\nSet-CellText $tbl4 3 5 \"46.1\"       # Stirrup spacing along width - Max.\nSet-CellText $tbl4 4 3 \"2.83\"       # Minimum shear reinforcement - Value\nSet-CellText $tbl4 4 4 \"1.75\"       # Minimum shear reinforcement - Min.\n\n$tbl5 = $d.Tables.Item(5)   # Shear reinforcement strength\nSet-CellText $tbl5 3 3 \"6.0\"        # Stirrup diameter (db)\nSet-CellText $tbl5 4 3 \"20.0\"       # Stirrup spacing (s)\nSet-CellText $tbl5 5 3 \"46.1\"       # Effective height (d)\nSet-CellText $tbl5 6 3 \"1.75\"       # Minimum shear reinforcing (Av,min)\nSet-CellText $tbl5 7 3 \"2.98\"       # Required shear reinforcing (Av,req)\nSet-CellText $tbl5 8 1 \"Defined shear reinforcing\"  # label\nSet-CellText $tbl5 8 3 \"2.83\"       # Shear reinforcing (Av)\nSet-CellText $tbl5 9 3 \"39.1\"       # Shear steel strength (\u00d8Vs)\n\n$tbl6 = $d.Tables.Item(6)   # Concrete strength\nSet-CellText $tbl6 2 3 \"922.0\"      # Effective shear area (Acv)\nSet-CellText $tbl6 3 3 \"0.0\"        # Longitudinal reinforcement ratio (\u03c1w)\nSet-CellText $tbl6 4 3 \"0.839\"      # Size modification factor (\u03bbs)\nSet-CellText $tbl6 5 3 \"0.0\"        # Axial stress (\u03c3Nu)\nSet-CellText $tbl6 6 3 \"0.85\"       # Concrete effective shear stress (kc)\nSet-CellText $tbl6 7 3 \"58.78\"      # Concrete strength (\u00d8Vc)\nSet-CellText $tbl6 8 3 \"286.97\"     # Maximum shear strength (\u00d8Vmax)\nSet-CellText $tbl6 9 3 \"97.88\"      # Total shear strength (\u00d8Vn)\nSet-CellText $tbl6 11 3 \"1.02\"      # Demand Capacity Ratio (DCR)\nSet-CellText $tbl6 11 4 \"\u274c\"        # Demand Capacity Ratio - Ok?\n"}
